# Bump the "Förändrad" (changed) date in column C for every data row
# (rows 2-61) from 45180 (2023-09-11) to 45181 (2023-09-12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C61").Value = 45181
